$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.849.25'
$ws.Range('E2').Value = '  +1.66%  '

$ws.Range('D3').Value = '1.667.55'
$ws.Range('E3').Value = '  +0.76%  '

$ws.Range('D4').Value = "'1.005"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.24%  '

$ws.Range('D5').Value = "'329.56"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +7.08%  '

$ws.Range('D6').Value = "'1.002"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.26%  '

$ws.Range('D7').Value = "'0.3642"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.43%  '

$ws.Range('D8').Value = "'46.69"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.35%  '

$ws.Range('D9').Value = "'0.3230"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.12%  '

$ws.Range('D10').Value = "'1.138"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.09%  '

$ws.Range('D11').Value = "'0.07040"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.23%  '

$ws.Range('D12').Value = "'1.003"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.43%  '

$ws.Range('D13').Value = "'6.051"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.83%  '

$ws.Range('D14').Value = "'19.50"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.81%  '

$ws.Range('D15').Value = '1.662.90'
$ws.Range('E15').Value = '  +0.69%  '

$ws.Range('D16').Value = "'6.593"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.45%  '

$ws.Range('D17').Value = "'0.00001042"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.06%  '

$ws.Range('D18').Value = "'0.06552"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.57%  '

$ws.Range('D19').Value = "'1.002"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.29%  '

$ws.Range('D20').Value = "'78.50"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.68%  '

$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').Value = "'15.77"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.25%  '

$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = "'5.899"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.43%  '

$ws.Range('D23').Value = "'12.93"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.14%  '

$ws.Range('D24').Value = '24.867.21'
$ws.Range('E24').Value = '  +1.96%  '

$ws.Range('D25').Value = "'2.438"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.41%  '

$ws.Range('D26').Value = "'2.377"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.40%  '

$ws.Range('D27').Value = "'148.10"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.37%  '

$ws.Range('D28').Value = "'18.62"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.13%  '

$ws.Range('D29').Value = '1.848.17'
$ws.Range('E29').Value = '  +0.60%  '

$ws.Range('D30').Value = "'125.43"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.12%  '

$ws.Range('D31').Value = "'1.172"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.24%  '

$ws.Range('D32').Value = "'4.079"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.70%  '

$ws.Range('D33').Value = "'5.713"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.67%  '

$ws.Range('D34').Value = "'0.08421"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.14%  '

$ws.Range('D35').Value = "'1.640"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.01%  '

$ws.Range('D36').Value = "'12.19"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.21%  '

$ws.Range('D37').Value = "'5.130"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.56%  '

$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').Value = "'1.229"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.96%  '

$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = "'0.06009"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.76%  '

$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = "'0.02223"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.90%  '

$ws.Range('D41').Value = "'0.2079"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.51%  '

$ws.Range('D42').Value = "'8.182"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.20%  '

$ws.Range('D43').Value = "'1.001"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.26%  '

$ws.Range('D44').Value = "'0.5928"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.33%  '

$ws.Range('D45').Value = "'13.67"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +8.66%  '

$ws.Range('D46').Value = "'3.846"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.02%  '

$ws.Range('D47').Value = "'0.5713"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.42%  '

$ws.Range('D48').Value = "'124.55"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.91%  '

$ws.Range('D49').Value = "'1.954"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.84%  '

$ws.Range('D50').Value = "'0.06984"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.16%  '

$ws.Range('D51').Value = "'1.185"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.62%  '
